$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.711.36'
$ws.Range("E2").Value = '  -2.73%  '

$ws.Range("D3").Value = '2.092.79'
$ws.Range("E3").Value = '  -1.98%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '344.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5156'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.97%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4366'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.97'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09223'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.161'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.87%  '

$ws.Range("D13").Value = '2.098.25'
$ws.Range("E13").Value = '  -1.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.262'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.734'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001148'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.011'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06652'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.93%  '

$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.177'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.03%  '

$ws.Range("D23").Value = '29.732.76'
$ws.Range("E23").Value = '  -3.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.320'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.12%  '

$ws.Range("D26").Value = '2.343.65'
$ws.Range("E26").Value = '  -2.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.511'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.76'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.50%  '

$ws.Range("E31").Value = '  -7.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1050'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.643'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.143'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.937'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.237'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02556'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06677'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.20%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.345'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.20%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.59%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6852'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.80%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.2221'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.75%  '

$ws.Range("E44").Value = '  +2.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.56%  '

$ws.Range("E46").Value = '  -2.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000359'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.623'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.43%  '

$ws.Range("E49").Value = '  -2.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '81.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.163'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.99%  '
